$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-7 down to 5-8
$ws.Rows.Item(4).Insert()

# Populate new row 4 with the new match data (Botev Plovdiv vs Cherno More)
$ws.Range("A4").Value = "hKMIYnqh"
$ws.Range("B4").Value = "25/10/2024"
$ws.Range("C4").Value = "10:30"
$ws.Range("D4").Value = "BULGARIA - PARVA LIGA"
$ws.Range("E4").Value = "Botev Plovdiv"
$ws.Range("F4").Value = "Cherno More"
$ws.Range("G4").Value = 2.4
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 3.3
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 1.91
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 10
$ws.Range("Y4").Value = 10
$ws.Range("Z4").Value = 23
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 5.5
$ws.Range("AE4").Value = 17
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 351
$ws.Range("AH4").Value = 8
$ws.Range("AI4").Value = 15
$ws.Range("AJ4").Value = 13
$ws.Range("AK4").Value = 34
$ws.Range("AL4").Value = 34
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 4.33
$ws.Range("AO4").Value = 15
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 81
$ws.Range("AS4").Value = 301
$ws.Range("AT4").Value = 2.25
$ws.Range("AU4").Value = 9
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 51
$ws.Range("AX4").Value = 5
$ws.Range("AY4").Value = 19
$ws.Range("AZ4").Value = 34
$ws.Range("BA4").Value = 67
$ws.Range("BB4").Value = 101
$ws.Range("BC4").Value = 301
$ws.Range("BD4").Value = 51

# Update odds that changed for the Penang vs Sabah match (now row 5)
$ws.Range("G5").Value = 3.15
$ws.Range("H5").Value = 3.25
$ws.Range("J5").Value = 3.55
$ws.Range("K5").Value = 2.18
$ws.Range("O5").Value = 1.17
$ws.Range("P5").Value = 4.32
$ws.Range("T5").Value = 3.27
$ws.Range("W5").Value = 10.5
$ws.Range("X5").Value = 16
$ws.Range("Y5").Value = 9.25
$ws.Range("Z5").Value = 35
$ws.Range("AA5").Value = 20
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 5.9
$ws.Range("AE5").Value = 9.5
$ws.Range("AF5").Value = 30
$ws.Range("AG5").Value = 150
$ws.Range("AI5").Value = 10
$ws.Range("AJ5").Value = 7.2
$ws.Range("AK5").Value = 17
$ws.Range("AM5").Value = 16
$ws.Range("AN5").Value = 5.4
$ws.Range("AO5").Value = 17
$ws.Range("AQ5").Value = 75
$ws.Range("AR5").Value = 90
$ws.Range("AT5").Value = 3.1
$ws.Range("AU5").Value = 6.2
$ws.Range("AV5").Value = 40
$ws.Range("AX5").Value = 4.25
$ws.Range("AY5").Value = 10.5
$ws.Range("AZ5").Value = 15.5
$ws.Range("BA5").Value = 37
$ws.Range("BC5").Value = 150

# Update the odds that changed for the Rukh Lviv vs Ch. Odesa match (now row 8)
$ws.Range("N8").Value = 7
